$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Activate()

$ws.Range("B17").Value = "Penalties, Fees, Interest, Principal order"
$ws.Range("B17").HorizontalAlignment = -4131  # xlLeft
$ws.Range("B17").VerticalAlignment = -4160    # xlTop

$ws.Range("B17").Select()
